# Sync GTD brain 2026-02-18
# Applies the edits described by the commit diff to the TBH Dispatch / Billing
# Verification Pack workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Cover sheet: bump the "Generated" timestamp
# ---------------------------------------------------------------------------
$cover = $wb.Worksheets.Item("Cover")
$cover.Range("B4").Value = "2026-02-18 09:43"

# ---------------------------------------------------------------------------
# 2. "Dispatch Plant Day" sheet: zero out column D (Quantity) for every row
#    that currently has a non-zero quantity.
# ---------------------------------------------------------------------------
$day = $wb.Worksheets.Item("Dispatch Plant Day")
$dayRows = @(3,4,5,7,8,9,10,14,15,16,17,18,19,20,21,23,24,25,26,27,28,29,30,
    32,33,34,35,37,38,39,40,43,44,45,46,48,49,50,51,54,57,59,60,61,62,64,65,
    66,67,70,71,72,73,76,77,78,79,80,81,82,83,86,88,89,91)
foreach ($r in $dayRows) {
    $day.Cells.Item($r, 4).Value = 0
}

# ---------------------------------------------------------------------------
# 3. "Dispatch Plant Month" sheet: zero out column F (Concrete Delivered Qty)
#    for the monthly roll-up rows.
# ---------------------------------------------------------------------------
$month = $wb.Worksheets.Item("Dispatch Plant Month")
$monthRows = @(2,3,4,9)
foreach ($r in $monthRows) {
    $month.Cells.Item($r, 6).Value = 0
}

# ---------------------------------------------------------------------------
# 4. "Dispatch vs AR (Invoice)" sheet
# ---------------------------------------------------------------------------
$ar = $wb.Worksheets.Item("Dispatch vs AR (Invoice)")

# 4a. Correct an existing invoice (row 630 / invoice 1483257): dispatch
#     revenue was restated, which changes the Difference column too.
$ar.Cells.Item(630, 2).Value = 3858.75
$ar.Cells.Item(630, 4).Value = -311.18

# 4b. Append newly-surfaced invoice lines at the bottom of the table
#     (rows 720-725). Column A holds the invoice code as text (it must
#     match the existing "number-looking text" storage used by the rest
#     of column A), columns B-D are currency values.
$newInvoices = @(
    @{ Row = 720; Code = "1483667"; Dispatch = 2970;     AR = 0; Diff = 2970 },
    @{ Row = 721; Code = "1483668"; Dispatch = -3858.75; AR = 0; Diff = -3858.75 },
    @{ Row = 722; Code = "1484728"; Dispatch = 535.5;    AR = 0; Diff = 535.5 },
    @{ Row = 723; Code = "1484729"; Dispatch = -535.5;   AR = 0; Diff = -535.5 },
    @{ Row = 724; Code = "1484730"; Dispatch = 1695.75;  AR = 0; Diff = 1695.75 },
    @{ Row = 725; Code = "1484731"; Dispatch = -1695.75; AR = 0; Diff = -1695.75 }
)

foreach ($inv in $newInvoices) {
    $r = $inv.Row

    $codeCell = $ar.Cells.Item($r, 1)
    $codeCell.NumberFormat = '@'
    $codeCell.Value = $inv.Code
    $codeCell.Style = 'Normal'

    $bCell = $ar.Cells.Item($r, 2)
    $bCell.NumberFormat = '$#,##0.00'
    $bCell.Value = $inv.Dispatch

    $cCell = $ar.Cells.Item($r, 3)
    $cCell.NumberFormat = '$#,##0.00'
    $cCell.Value = $inv.AR

    $dCell = $ar.Cells.Item($r, 4)
    $dCell.NumberFormat = '$#,##0.00'
    $dCell.Value = $inv.Diff
}
